# "Created the Special Statistics section" - trims today's order log down to
# the two most recent orders (new data) and drops the older rows entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the six now-obsolete data rows (old rows 4-9) first, so the sheet's
# used range / dimension shrinks to match the two remaining orders.
$ws.Range("A4:I9").EntireRow.Delete()

# Refresh row 2 with the newest order.
$ws.Range("A2").Value = 650
$ws.Range("B2").Value = 45761.22928240741
$ws.Range("C2").Value = "Ajay Francis Anchan"
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 400
$ws.Range("F2").Value = 10
$ws.Range("G2").Value = 10
$ws.Range("H2").Value = 20
$ws.Range("I2").Value = "Strawberry Lassi (x1), Butterscotch Lassi (x6)"

# Refresh row 3 with the next order.
$ws.Range("A3").Value = 649
$ws.Range("B3").Value = 45761.22928240741
$ws.Range("C3").Value = "Ajay Francis Anchan"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 350
$ws.Range("F3").Value = 2.3
$ws.Range("G3").Value = 2.3
$ws.Range("H3").Value = 10
$ws.Range("I3").Value = "Chicken Cheese Burger (x3), Chicken Wrap (x2), Banana Shake (x1)"
